$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "Objetivos:" value (row 10, cols B/C) ---
$objetivosPt = "1) Formativos: Propiciar ao educando as condições básicas e necessárias para a sua formação profissional. 2) Informativos: fornecer ao educando os conceitos básicos para o entendimento, assessoramento e acompanhamento de Projetos na Indústria Química seguindo metodologia especifica.3) Automatizantes: desenvolver no educando o raciocínio analítico, obedecendo metodologia sistemática aplicada em projetos."
$ws.Range("B10").Value2 = $objetivosPt
$ws.Range("C10").Value2 = $objetivosPt

# --- 2. Insert a new row at 13 (pushes old rows 13-24 down to 14-25) ---
$ws.Rows.Item(13).Insert()

# Remove the leftover formatted-but-empty cell Excel carries into column A of the new row
$ws.Range("A13").Clear()

# Populate the new row 13 (B/C) with the "Docentes responsáveis" name, matching the
# look (style) of the other text cells in columns B/C
$ws.Range("B2").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$docente = "5840671 - Francisco José Moreira Chaves"
$ws.Range("B13").Value2 = $docente
$ws.Range("C13").Value2 = $docente

# --- 3. "Programa resumido:" value -> now at row 14 after the shift ---
$programaResumido = "Projetos: Conceituação e Viabilidade, Ponto Nivelamento, Legalização Industrial."
$ws.Range("B14").Value2 = $programaResumido
$ws.Range("C14").Value2 = $programaResumido

# --- 4. "Programa:" value -> now at row 16 after the shift ---
$programa = "1.Introdução: Conceitos de Gestão de Projetos2.Aspectos da Implantação de Projetos: Etapas Fundamentais e Formas Parciais3.O Ciclo de Vida do Projeto4.Aspectos da Viabilidade de Projetos: Receitas, Custos, Ponto de Nivelamento, Estimativas5.Guia PMBOK: Principais Áreas de Conhecimento6.Plano de Projeto7.O Gerente de Projeto e as Interfaces: Equipes de Projeto8.Legalização de Projetos: Aspectos sobre o Licenciamento9.Técnicas de Análise de Riscos Operacionais10.Aspectos sobre Auditorias e Auditorias Integradas11.Apresentação de Trabalhos."
$ws.Range("B16").Value2 = $programa
$ws.Range("C16").Value2 = $programa

# --- 5. "Método:" value -> now at row 19 after the shift ---
$metodo = "Por meio de aulas presenciais, com apresentação dos fundamentos e exemplos ou casos, e também pela apresentação de trabalhos em equipes.Justificativa: adequação do método de avaliação."
$ws.Range("B19").Value2 = $metodo
$ws.Range("C19").Value2 = $metodo

# --- 6. "Critério:" value -> now at row 20 after the shift ---
$criterio = "A Avaliação será: MF = (P1 + P2)/2; Onde: P1: Trabalho; P2: Trabalho. Poderá haver também prova individual sobre os fundamentos."
$ws.Range("B20").Value2 = $criterio
$ws.Range("C20").Value2 = $criterio

# --- 7. "Norma de recuperação:" value -> now at row 21 after the shift ---
$normaRecuperacao = "Prova de exame."
$ws.Range("B21").Value2 = $normaRecuperacao
$ws.Range("C21").Value2 = $normaRecuperacao

# --- 8. "Bibliografia:" value -> now at row 22 after the shift ---
$bibliografia = @"
1) Chemical Engineering Plant Design (Vilbrandt e Bryden)
2) Project Engineering of Process Plants Rose e Barrow (2ª impressão - 1968)
3) Elaboração e Análise de Projetos Simonsen, M.H. e H. Flanger
4) Implantação de Indústrias. Valle, E.C. Livros Técnicos e Científicos Editores S/A, Rio de Janeiro. 
5) Introdução ao Projeto de Engenharia. Asimov, Morris. Editora Mestre Jou-São Paulo, 1968.
6) IBGR  Instituto Brasileiro de Gerenciamento de Riscos, 2000.
"@
$bibliografia = $bibliografia.TrimEnd("`r", "`n")
$ws.Range("B22").Value2 = $bibliografia
$ws.Range("C22").Value2 = $bibliografia
